$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# dataset_internal_id: LandAndGender -> LG
$ws.Range("B2").Value = "LG"

# indicator_internal_id: LandAndGender.2M -> LG.2M
$ws.Range("B3").Value = "LG.2M"
